$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Holidays 2019")

# Update Date_of_sales (column I) for rows 2 through 21 from 45110 (2023-07-03)
# to 45138 (2023-07-31), keeping the existing date number format.
$ws.Range("I2:I21").Value = 45138
